# ---------------------------------------------------------------------------
# Commit "version2 with changes":
#  - adds five new Python data-structure sheets (Tree/Graph/Stack/Queue/
#    LinkedList), modelled on the existing ArraypythonCode sheet
#  - leaves Login/ArraypythonCode content untouched but updates their saved
#    cursor/selection
#  - ends with GraphpythonCode active (tab scrolled so later tabs show)
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$xlPasteFormats  = -4122
$xlVAlignBottom  = -4107

$login = $wb.Worksheets.Item("Login")
$arr   = $wb.Worksheets.Item("ArraypythonCode")

# ---------------------------------------------------------------------------
# 1. Create the five new worksheets, in order, after ArraypythonCode.
# ---------------------------------------------------------------------------
$newNames = @("TreepythonCode","GraphpythonCode","StackpythonCode","Queuecode","Linkedlistcode")
foreach ($name in $newNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $created = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $created.Name = $name
}

$tree   = $wb.Worksheets.Item("TreepythonCode")
$graph  = $wb.Worksheets.Item("GraphpythonCode")
$stack  = $wb.Worksheets.Item("StackpythonCode")
$queue  = $wb.Worksheets.Item("Queuecode")
$linked = $wb.Worksheets.Item("Linkedlistcode")

# ---------------------------------------------------------------------------
# Helper content shared by every new sheet: row3/A3 always repeats the
# "print('Hello, world!)" string already used on ArraypythonCode.
# ---------------------------------------------------------------------------
$helloText = "print('Hello, world!)"

# ---------------------------------------------------------------------------
# 2. TreepythonCode  (A1:C6)
# ---------------------------------------------------------------------------
$tree.Range("A1").Value = "TreepythonCode"
$arr.Range("A1").Copy()
$tree.Range("A1").PasteSpecial($xlPasteFormats)

$tree.Range("A2").Value = 'print("This is Tree in python")'
$arr.Range("A2").Copy()
$tree.Range("A2").PasteSpecial($xlPasteFormats)
$tree.Rows.Item(2).RowHeight = 28

$tree.Range("A3").Value = $helloText

# stray formatted-but-empty cell, mirrors the wrap+vcenter style used elsewhere
$tree.Range("C6").VerticalAlignment = $xlVAlignBottom
$tree.Range("C6").WrapText = $true
$tree.Range("C6").VerticalAlignment = -4108
$excel.CutCopyMode = $false

$tree.Columns.Item(1).ColumnWidth = 16.666666666666668

# ---------------------------------------------------------------------------
# 3. GraphpythonCode  (A1:A3)
# ---------------------------------------------------------------------------
$graph.Range("A1").Value = "GraphpythonCode"

$graph.Range("A2").Value = 'print("This is Graph in python")'
$arr.Range("A2").Copy()
$graph.Range("A2").PasteSpecial($xlPasteFormats)
$graph.Rows.Item(2).RowHeight = 28

$graph.Range("A3").Value = $helloText
$excel.CutCopyMode = $false

$graph.Columns.Item(1).ColumnWidth = 15.916666666666666

# ---------------------------------------------------------------------------
# 4. StackpythonCode  (A1:A3)
# ---------------------------------------------------------------------------
$stack.Range("A1").Value = "StackpythonCode"

$stack.Range("A2").Value = 'print("This is Stack in python")'
$arr.Range("A2").Copy()
$stack.Range("A2").PasteSpecial($xlPasteFormats)
$stack.Rows.Item(2).RowHeight = 28

$stack.Range("A3").Value = $helloText
$excel.CutCopyMode = $false

$stack.Columns.Item(1).ColumnWidth = 15.541666666666666

# ---------------------------------------------------------------------------
# 5. Queuecode  (A1:A3)
# ---------------------------------------------------------------------------
$queue.Range("A1").Value = "Queuecode"

$queue.Range("A2").Value = 'print("This is Queue in python")'
$arr.Range("A2").Copy()
$queue.Range("A2").PasteSpecial($xlPasteFormats)
$queue.Rows.Item(2).RowHeight = 28

$queue.Range("A3").Value = $helloText
$excel.CutCopyMode = $false

$queue.Columns.Item(1).ColumnWidth = 15.916666666666666

# ---------------------------------------------------------------------------
# 6. Linkedlistcode  (A1:A3)
# ---------------------------------------------------------------------------
$linked.Range("A1").Value = "Linkedlistcode"

$linked.Range("A2").Value = 'print("This is LinkedList in python")'
$arr.Range("A2").Copy()
$linked.Range("A2").PasteSpecial($xlPasteFormats)
$linked.Rows.Item(2).RowHeight = 42

$linked.Range("A3").Value = $helloText
$excel.CutCopyMode = $false

$linked.Columns.Item(1).ColumnWidth = 15.916666666666666

# ---------------------------------------------------------------------------
# 7. Restore / update cursor positions on every sheet, finishing on
#    GraphpythonCode!F18 so it ends up the active tab (matches tabSelected
#    moving off Login and onto GraphpythonCode).
# ---------------------------------------------------------------------------
$login.Range("C15").Select()
$arr.Range("A2").Select()
$tree.Range("A2:A3").Select()
$stack.Range("A2:A3").Select()
$queue.Range("A2:A3").Select()
$linked.Range("E6").Select()
$graph.Range("F18").Select()
